$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (the default/no-explicit-style used by existing data rows)
$normalStyle = $ws.Range("A3").Style

# ---- Row 2 ----
$ws.Range("F2").Value = 1298000
$ws.Range("G2").Value = 1298000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 10400
$ws.Range("J2").Value = 12700
$ws.Range("K2").Value = 8650735
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 16000
$ws.Range("P2").Value = 2055746777
$ws.Range("Q2").Value = 2918221978
$ws.Range("R2").Value = 979078233
$ws.Range("S2").Value = -3343774083
$ws.Range("T2").Value = -3525649863
$ws.Range("U2").Value = -1713494359
$ws.Range("V2").Value = -4430074915
$ws.Range("W2").Value = -8304699942
$ws.Range("X2").Value = -1627684107
foreach ($addr in @("A2","B2","C2","D2","E2","N2","O2","Y2")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A2").Value = '2024-05-27'
$ws.Range("B2").Value = '2024-05-31'
$ws.Range("C2").Value = '2024-06-17'
$ws.Range("D2").Value = '대신'
$ws.Range("E2").Value = '라메디텍'
$ws.Range("N2").Value = '1115.44:1'
$ws.Range("O2").Value = '9.93%'
$ws.Range("Y2").Value = '초소형 레이저 의료기기 및 미용기기'
foreach ($addr in @("A2","B2","C2","D2","E2","N2","O2","Y2")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 3 ----
$ws.Range("F3").Value = 1400000
$ws.Range("G3").Value = 1400000
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 34000
$ws.Range("J3").Value = 40000
$ws.Range("K3").Value = 7942750
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 40000
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
foreach ($addr in @("A3","B3","C3","D3","E3","N3","O3","Y3")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A3").Value = '2024-05-23'
$ws.Range("B3").Value = '2024-05-29'
$ws.Range("C3").Value = '2024-06-14'
$ws.Range("D3").Value = '삼성'
$ws.Range("E3").Value = '그리드위즈'
$ws.Range("N3").Value = '124.60:1'
$ws.Range("O3").Value = '0.95%'
$ws.Range("Y3").Value = '수요관리 서비스, 전기차 충전기 모뎀 등'
foreach ($addr in @("A3","B3","C3","D3","E3","N3","O3","Y3")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 4 ----
$ws.Range("F4").Value = 6650000
$ws.Range("G4").Value = 6650000
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 8100000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 2000
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
foreach ($addr in @("A4","B4","C4","D4","E4","N4","O4","Y4")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A4").Value = '2024-05-13'
$ws.Range("B4").Value = '2024-05-14'
$ws.Range("C4").Value = '2024-05-29'
$ws.Range("D4").Value = '미래'
$ws.Range("E4").Value = '미래에셋비전스팩4호'
$ws.Range("N4").Value = '1011.2:1'
$ws.Range("O4").Value = '-'
$ws.Range("Y4").Value = '기업인수목적회사(기타금융서비스)'
foreach ($addr in @("A4","B4","C4","D4","E4","N4","O4","Y4")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 5 ----
$ws.Range("F5").Value = 1200000
$ws.Range("G5").Value = 1200000
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 8700
$ws.Range("J5").Value = 11000
$ws.Range("K5").Value = 7651263
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 14000
$ws.Range("P5").Value = 468321534076
$ws.Range("Q5").Value = 555936831337
$ws.Range("R5").Value = 359249623614
$ws.Range("S5").Value = 22403886436
$ws.Range("T5").Value = 33386727728
$ws.Range("U5").Value = 10411712773
$ws.Range("V5").Value = 10859975142
$ws.Range("W5").Value = 29346086803
$ws.Range("X5").Value = 4820429371
foreach ($addr in @("A5","B5","C5","D5","E5","N5","O5","Y5")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A5").Value = '2024-04-30'
$ws.Range("B5").Value = '2024-05-08'
$ws.Range("C5").Value = '2024-05-23'
$ws.Range("D5").Value = '삼성'
$ws.Range("E5").Value = '노브랜드'
$ws.Range("N5").Value = '1075.61:1'
$ws.Range("O5").Value = '4.51%'
$ws.Range("Y5").Value = 'Knit, Woven 의류'
foreach ($addr in @("A5","B5","C5","D5","E5","N5","O5","Y5")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 6 ----
$ws.Range("F6").Value = 5000000
$ws.Range("G6").Value = 5000000
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 5505000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0
foreach ($addr in @("A6","B6","C6","D6","E6","N6","O6","Y6")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A6").Value = '2024-04-29'
$ws.Range("B6").Value = '2024-04-30'
$ws.Range("C6").Value = '2024-05-17'
$ws.Range("D6").Value = 'KB'
$ws.Range("E6").Value = 'KB제28호스팩'
$ws.Range("N6").Value = '1118.39:1'
$ws.Range("O6").Value = '-'
$ws.Range("Y6").Value = '기업인수합병'
foreach ($addr in @("A6","B6","C6","D6","E6","N6","O6","Y6")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 7 ----
$ws.Range("F7").Value = 1970000
$ws.Range("G7").Value = 1970000
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 13000
$ws.Range("J7").Value = 16000
$ws.Range("K7").Value = 13124496
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0
foreach ($addr in @("A7","B7","C7","D7","E7","N7","O7","Y7")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A7").Value = '2024-04-24'
$ws.Range("B7").Value = '2024-05-30'
$ws.Range("C7").Value = '2024-05-17'
$ws.Range("D7").Value = 'NH'
$ws.Range("E7").Value = '아이씨티케이'
$ws.Range("N7").Value = '783.2:1'
$ws.Range("O7").Value = '6.54%'
$ws.Range("Y7").Value = 'PUF반도체,보안솔루션(보안반도체,정보통신모듈기기,정보통신용반도체) 제조,개발'
foreach ($addr in @("A7","B7","C7","D7","E7","N7","O7","Y7")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 8 ----
$ws.Range("F8").Value = 1500000
$ws.Range("G8").Value = 1500000
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 11000
$ws.Range("J8").Value = 14000
$ws.Range("K8").Value = 8503460
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 18000
$ws.Range("P8").Value = 47284698907
$ws.Range("Q8").Value = 38750429966
$ws.Range("R8").Value = 25900014771
$ws.Range("S8").Value = 7595091433
$ws.Range("T8").Value = 5807002440
$ws.Range("U8").Value = 3668321605
$ws.Range("V8").Value = 5701880294
$ws.Range("W8").Value = 4780312126
$ws.Range("X8").Value = 4195570793
foreach ($addr in @("A8","B8","C8","D8","E8","N8","O8","Y8")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A8").Value = '2024-04-15'
$ws.Range("B8").Value = '2024-04-19'
$ws.Range("C8").Value = '2024-05-07'
$ws.Range("D8").Value = '한국'
$ws.Range("E8").Value = '코칩'
$ws.Range("N8").Value = '988.32:1'
$ws.Range("O8").Value = '13.19%'
$ws.Range("Y8").Value = '소형 및 초소형 슈퍼커패시터'
foreach ($addr in @("A8","B8","C8","D8","E8","N8","O8","Y8")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 9 ----
$ws.Range("F9").Value = 3000000
$ws.Range("G9").Value = 3000000
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 2000
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 3310000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 2000
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = 0
foreach ($addr in @("A9","B9","C9","D9","E9","N9","O9","Y9")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A9").Value = '2024-04-17'
$ws.Range("B9").Value = '2024-04-18'
$ws.Range("C9").Value = '2024-05-07'
$ws.Range("D9").Value = 'SK'
$ws.Range("E9").Value = 'SK증권제12호스팩'
$ws.Range("N9").Value = '1,189.41:1'
$ws.Range("O9").Value = '-'
$ws.Range("Y9").Value = '기업인수목적 주식회사'
foreach ($addr in @("A9","B9","C9","D9","E9","N9","O9","Y9")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 10 ----
$ws.Range("F10").Value = 3000000
$ws.Range("G10").Value = 3000000
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 6500
$ws.Range("J10").Value = 8500
$ws.Range("K10").Value = 21945300
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 10500
$ws.Range("P10").Value = 9576212189
$ws.Range("Q10").Value = 11914994171
$ws.Range("R10").Value = 7285537916
$ws.Range("S10").Value = 1172310325
$ws.Range("T10").Value = -2762203259
$ws.Range("U10").Value = -4737405164
$ws.Range("V10").Value = -7460336546
$ws.Range("W10").Value = -7104430732
$ws.Range("X10").Value = -7501425172
foreach ($addr in @("A10","B10","C10","D10","E10","N10","O10","Y10")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A10").Value = '2024-04-12'
$ws.Range("B10").Value = '2024-04-18'
$ws.Range("C10").Value = '2024-05-03'
$ws.Range("D10").Value = 'KB'
$ws.Range("E10").Value = '민테크'
$ws.Range("N10").Value = '946.72:1'
$ws.Range("O10").Value = '4.23%'
$ws.Range("Y10").Value = '배터리 진단시스템, 배터리 시스템, 충방전 검사장비'
foreach ($addr in @("A10","B10","C10","D10","E10","N10","O10","Y10")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 11 ----
$ws.Range("F11").Value = 1100000
$ws.Range("G11").Value = 1100000
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 22000
$ws.Range("J11").Value = 26000
$ws.Range("K11").Value = 10429232
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 33000
$ws.Range("P11").Value = -75676750274
$ws.Range("Q11").Value = -68652978862
$ws.Range("R11").Value = -9506668082
$ws.Range("S11").Value = -69862474811
$ws.Range("T11").Value = -137025491259
$ws.Range("U11").Value = 3014576074
$ws.Range("V11").Value = 0
$ws.Range("W11").Value = 0
$ws.Range("X11").Value = 0
foreach ($addr in @("A11","B11","C11","D11","E11","N11","O11","Y11")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A11").Value = '2024-04-12'
$ws.Range("B11").Value = '2024-04-18'
$ws.Range("C11").Value = '2024-05-02'
$ws.Range("D11").Value = '한국'
$ws.Range("E11").Value = '디앤디파마텍'
$ws.Range("N11").Value = '848.50:1'
$ws.Range("O11").Value = '10.96%'
$ws.Range("Y11").Value = '대사성질환 치료제 등'
foreach ($addr in @("A11","B11","C11","D11","E11","N11","O11","Y11")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 12 ----
$ws.Range("F12").Value = 5150000
$ws.Range("G12").Value = 5150000
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 5510000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 2000
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("T12").Value = 0
$ws.Range("U12").Value = 0
$ws.Range("V12").Value = 0
$ws.Range("W12").Value = 0
$ws.Range("X12").Value = 0
foreach ($addr in @("A12","B12","C12","D12","E12","N12","O12","Y12")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A12").Value = '2024-04-15'
$ws.Range("B12").Value = '2024-04-16'
$ws.Range("C12").Value = '2024-05-02'
$ws.Range("D12").Value = '유안타'
$ws.Range("E12").Value = '유안타제16호스팩'
$ws.Range("N12").Value = '1,050.42:1'
$ws.Range("O12").Value = '-'
$ws.Range("Y12").Value = '금융 지원 서비스(기업인수목적회사)'
foreach ($addr in @("A12","B12","C12","D12","E12","N12","O12","Y12")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 13 ----
$ws.Range("F13").Value = 3500000
$ws.Range("G13").Value = 3500000
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 3700000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 2000
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("U13").Value = 0
$ws.Range("V13").Value = 0
$ws.Range("W13").Value = 0
$ws.Range("X13").Value = 0
foreach ($addr in @("A13","B13","C13","D13","E13","N13","O13","Y13")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A13").Value = '2024-04-08'
$ws.Range("B13").Value = '2024-04-09'
$ws.Range("C13").Value = '2024-04-24'
$ws.Range("D13").Value = '하나'
$ws.Range("E13").Value = '하나33호스팩'
$ws.Range("N13").Value = '1277.22:1'
$ws.Range("O13").Value = '-'
$ws.Range("Y13").Value = '기업인수합병'
foreach ($addr in @("A13","B13","C13","D13","E13","N13","O13","Y13")) { $ws.Range($addr).Style = $normalStyle }

# ---- Row 14 ----
$ws.Range("F14").Value = 3000000
$ws.Range("G14").Value = 3000000
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 3620000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 2000
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 0
$ws.Range("U14").Value = 0
$ws.Range("V14").Value = 0
$ws.Range("W14").Value = 0
$ws.Range("X14").Value = 0
foreach ($addr in @("A14","B14","C14","D14","E14","N14","O14","Y14")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A14").Value = '2024-04-04'
$ws.Range("B14").Value = '2024-04-05'
$ws.Range("C14").Value = '2024-04-22'
$ws.Range("D14").Value = '신한'
$ws.Range("E14").Value = '신한제13호스팩'
$ws.Range("N14").Value = '1337.88:1'
$ws.Range("O14").Value = '-'
$ws.Range("Y14").Value = '기타금융서비스(기업합병)'
foreach ($addr in @("A14","B14","C14","D14","E14","N14","O14","Y14")) { $ws.Range($addr).Style = $normalStyle }
